$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04763786555579896
$ws.Range("C2").Value = 0.002777888934908601
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 0.7047202576384011

$ws.Range("B3").Value = 0.003994804209775715
$ws.Range("C3").Value = 0.04240448674262143
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 12.60706245710986

$ws.Range("B4").Value = 0.6753301551942219
$ws.Range("C4").Value = 0.3127903958511391
$ws.Range("D4").Value = 337.1190423067083
$ws.Range("E4").Value = 8.660232485948974
$ws.Range("G4").Value = 346.7673953437026

$ws.Range("B5").Value = 1.459612070389937
$ws.Range("C5").Value = 1.667794583268128
$ws.Range("D5").Value = 3.900430680208489
$ws.Range("E5").Value = 0.496779210170732
$ws.Range("G5").Value = 7.524616544037286
